# Automatische test-sync: 2025-06-17 22:50:13
# Append the new mail-log entry (row 59) to the "Logs" sheet, extend the
# conditional-formatting ranges so they keep covering the data, and bump the
# "Overig" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$newRow = 59
$logs.Cells.Item($newRow, 1).Value = "Sollicitatie marketingfunctie"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Hierbij solliciteer ik voor de functie van marketeer. Zie bijlage voor CV."
$logs.Cells.Item($newRow, 4).Value = "Overig"
$logs.Cells.Item($newRow, 6).Value = "2025-06-17 22:49:15"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# Extend the two conditionalFormatting blocks (Categorie / Beantwoord) from
# row 58 to row 59 so the new row is covered too.
$catRange = $logs.Range("D2:D59")
$catConditions = $logs.Range("D2:D58").FormatConditions
$catConditions.Item(1).ModifyAppliesToRange($catRange)

$answeredRange = $logs.Range("G2:G59")
$answeredConditions = $logs.Range("G2:G58").FormatConditions
$answeredConditions.Item(1).ModifyAppliesToRange($answeredRange)

# Update the Dashboard summary count for the "Overig" category (16 -> 17).
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(3, 2).Value = 17
